$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on numeric-looking Price cells so they remain text (matches source inlineStr values)
foreach ($addr in @("D5", "D6", "D8", "D9", "D10", "D11", "D13", "D14", "D16", "D18", "D20", "D22", "D23", "D25", "D26", "D27", "D28", "D29", "D32", "D33", "D34", "D35", "D37", "D38", "D39", "D44", "D46", "D47", "D48", "D49", "D50", "D51")) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '35.627.18'
$ws.Range("E2").Value = '  -2.39%  '

$ws.Range("D3").Value = '1.982.65'
$ws.Range("E3").Value = '  -3.66%  '

$ws.Range("E4").Value = '  -0.14%  '

$ws.Range("D5").Value = '241.52'
$ws.Range("E5").Value = '  +0.19%  '

$ws.Range("D6").Value = '0.635'
$ws.Range("E6").Value = '  -3.50%  '

$ws.Range("E7").Value = '  -0.04%  '

$ws.Range("D8").Value = '56.51'
$ws.Range("E8").Value = '  +8.29%  '

$ws.Range("D9").Value = '59.94'
$ws.Range("E9").Value = '  +1.58%  '

$ws.Range("D10").Value = '0.359'
$ws.Range("E10").Value = '  -0.13%  '

$ws.Range("D11").Value = '0.0728'
$ws.Range("E11").Value = '  -2.32%  '

$ws.Range("E12").Value = '  -4.87%  '

$ws.Range("D13").Value = '0.901'
$ws.Range("E13").Value = '  -1.14%  '

$ws.Range("D14").Value = '14.21'
$ws.Range("E14").Value = '  -2.58%  '

$ws.Range("D15").Value = '2.265.17'
$ws.Range("E15").Value = '  -4.02%  '

$ws.Range("D16").Value = '5.23'
$ws.Range("E16").Value = '  -2.82%  '

$ws.Range("D17").Value = '1.985.16'
$ws.Range("E17").Value = '  -5.08%  '

$ws.Range("D18").Value = '17.18'
$ws.Range("E18").Value = '  +5.89%  '

$ws.Range("D19").Value = '35.443.97'
$ws.Range("E19").Value = '  -2.77%  '

$ws.Range("D20").Value = '70.15'
$ws.Range("E20").Value = '  -1.70%  '

$ws.Range("E21").Value = '  -2.24%  '

$ws.Range("D22").Value = '232.34'
$ws.Range("E22").Value = '  -1.53%  '

$ws.Range("D23").Value = '5.04'
$ws.Range("E23").Value = '  -3.41%  '

$ws.Range("E24").Value = '  +0.17%  '

$ws.Range("D25").Value = '2.27'
$ws.Range("E25").Value = '  -4.41%  '

$ws.Range("D26").Value = '2.26'
$ws.Range("E26").Value = '  +7.00%  '

$ws.Range("D27").Value = '163.28'
$ws.Range("E27").Value = '  -0.43%  '

$ws.Range("D28").Value = '9.08'
$ws.Range("E28").Value = '  -3.76%  '

$ws.Range("D29").Value = '19.45'
$ws.Range("E29").Value = '  -4.10%  '

$ws.Range("E30").Value = '  -2.73%  '

$ws.Range("E31").Value = '  -0.32%  '

$ws.Range("D32").Value = '4.78'
$ws.Range("E32").Value = '  -4.98%  '

$ws.Range("D33").Value = '0.0584'
$ws.Range("E33").Value = '  -1.18%  '

$ws.Range("D34").Value = '0.0892'
$ws.Range("E34").Value = '  +9.56%  '

$ws.Range("D35").Value = '4.26'
$ws.Range("E35").Value = '  -7.19%  '

$ws.Range("E36").Value = '  -0.19%  '

$ws.Range("D37").Value = '2.27'
$ws.Range("E37").Value = '  +0.29%  '

$ws.Range("D38").Value = '1.79'
$ws.Range("E38").Value = '  -2.65%  '

$ws.Range("D39").Value = '4.88'
$ws.Range("E39").Value = '  +1.83%  '

$ws.Range("E40").Value = '  -4.43%  '

$ws.Range("E41").Value = '  -3.16%  '

$ws.Range("E42").Value = '  -3.07%  '

$ws.Range("E43").Value = '  -4.56%  '

$ws.Range("D44").Value = '0.0890'
$ws.Range("E44").Value = '  -5.46%  '

$ws.Range("B45").Value = 'Maker'
$ws.Range("C45").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D45").Value = '1.377.15'
$ws.Range("E45").Value = '  -0.01%  '

$ws.Range("B46").Value = 'Aave'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D46").Value = '90.55'
$ws.Range("E46").Value = '  -2.77%  '

$ws.Range("D47").Value = '7.40'
$ws.Range("E47").Value = '  -0.08%  '

$ws.Range("D48").Value = '15.40'
$ws.Range("E48").Value = '  -0.13%  '

$ws.Range("D49").Value = '2.87'
$ws.Range("E49").Value = '  +0.37%  '

$ws.Range("D50").Value = '2.26'
$ws.Range("E50").Value = '  -3.66%  '

$ws.Range("D51").Value = '45.60'
$ws.Range("E51").Value = '  +3.13%  '
